# Apply the "identifying weak spots and improving matches" edit.
#
# Summary of the change (per the target diff):
#  - The four sheets (salinity, ocean_acidification, noise_pollution,
#    light_pollution) each have an "adcap_score" column (D) for the
#    Movement/R (Reproductive Traits) categories (rows 2-14 and 27-28)
#    that used to hold raw numeric scores (0, 0.25, 0.5, 0.75, 1). These
#    are now replaced with categorical labels "none"/"low"/"medium"/"high".
#  - For the "Spatial Scale of species" rows (16-26) the D column stays
#    numeric (0/1), with a couple of corrected values on some sheets, but
#    the repeated boilerplate note in column E
#    ("score these as 1 = ...") is removed from every row except the
#    first (row 16), and that remaining note is replaced by a
#    sheet-specific comment about the stressor.
#  - Sheet view selection/active-cell bookkeeping changes, and
#    light_pollution becomes the selected sheet/tab.

$wb = $excel.ActiveWorkbook

# ---- salinity ----
$ws1 = $wb.Worksheets.Item("salinity")
$scores1 = @{
    2  = "none"
    3  = "none"
    4  = "low"
    5  = "high"
    6  = "medium"
    7  = "high"
    8  = "high"
    9  = "high"
    10 = "medium"
    11 = "medium"
    12 = "low"
    13 = "none"
    14 = "none"
    27 = "none"
    28 = "high"
}
foreach ($row in $scores1.Keys) { $ws1.Range("D$row").Value = $scores1[$row] }
$ws1.Range("E16").Value = "assume salinity stressor is limited to coastal waters"
foreach ($row in 17..26) { $ws1.Range("E$row").Value = $null }
$ws1.Range("E32").Select()

# ---- ocean_acidification ----
$ws2 = $wb.Worksheets.Item("ocean_acidification")
$scores2 = @{
    2  = "none"
    3  = "none"
    4  = "none"
    5  = "medium"
    6  = "medium"
    7  = "high"
    8  = "high"
    9  = "high"
    10 = "medium"
    11 = "medium"
    12 = "low"
    13 = "none"
    14 = "none"
    27 = "none"
    28 = "high"
}
foreach ($row in $scores2.Keys) { $ws2.Range("D$row").Value = $scores2[$row] }
foreach ($row in 17..20) { $ws2.Range("D$row").Value = 1 }
$ws2.Range("E16").Value = "deep water already more acidic; current mixing across depths"
foreach ($row in 17..26) { $ws2.Range("E$row").Value = $null }
$ws2.Range("F24").Select()

# ---- noise_pollution ----
$ws3 = $wb.Worksheets.Item("noise_pollution")
$scores3 = @{
    2 = "none"
    3 = "none"
    4 = "none"
    5 = "medium"
    6 = "medium"
    7 = "high"
    8 = "high"
}
foreach ($row in $scores3.Keys) { $ws3.Range("D$row").Value = $scores3[$row] }
$ws3.Range("E16").Value = "perhaps shipping noise isn't important at depth but sonar pulses travel for miles"
foreach ($row in 17..26) { $ws3.Range("E$row").Value = $null }
$ws3.Range("E17:E26").Select()

# ---- light_pollution ----
$ws4 = $wb.Worksheets.Item("light_pollution")
$scores4 = @{
    2 = "none"
    3 = "none"
    4 = "low"
    5 = "medium"
    6 = "high"
    7 = "high"
    8 = "high"
}
foreach ($row in $scores4.Keys) { $ws4.Range("D$row").Value = $scores4[$row] }
$ws4.Range("E16").Value = "assume light pollution limited to the surface"
foreach ($row in 17..26) { $ws4.Range("E$row").Value = $null }
$ws4.Range("E32").Select()
$ws4.Activate()
